# "Generate Report for Archive"
#
# The localization-status report is regenerated: the two handoff rows that
# were still "Ready for handoff" have moved on to "In Translation", and the
# Status/zh-cn/de-de columns that hold that text are re-sized (narrower,
# since "In Translation" renders shorter than "Ready for handoff").

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: zh-cn (E) / de-de (F) status cells for both rows -----
if ($overview.Range("E2").Value2 -eq $oldStatus) { $overview.Range("E2").Value = $newStatus }
if ($overview.Range("F2").Value2 -eq $oldStatus) { $overview.Range("F2").Value = $newStatus }
if ($overview.Range("E3").Value2 -eq $oldStatus) { $overview.Range("E3").Value = $newStatus }
if ($overview.Range("F3").Value2 -eq $oldStatus) { $overview.Range("F3").Value = $newStatus }

# --- zh-cn / de-de sheets: Status column (C) for both rows ----------------
if ($zhcn.Range("C2").Value2 -eq $oldStatus) { $zhcn.Range("C2").Value = $newStatus }
if ($zhcn.Range("C3").Value2 -eq $oldStatus) { $zhcn.Range("C3").Value = $newStatus }

if ($dede.Range("C2").Value2 -eq $oldStatus) { $dede.Range("C2").Value = $newStatus }
if ($dede.Range("C3").Value2 -eq $oldStatus) { $dede.Range("C3").Value = $newStatus }

# --- Re-fit the columns that held the status text to their new, narrower
#     content width -----------------------------------------------------
$newWidth = 12.5

$overview.Columns.Item(5).ColumnWidth = $newWidth   # Overview!E (zh-cn)
$overview.Columns.Item(6).ColumnWidth = $newWidth   # Overview!F (de-de)
$zhcn.Columns.Item(3).ColumnWidth     = $newWidth   # zh-cn!C (Status)
$dede.Columns.Item(3).ColumnWidth     = $newWidth   # de-de!C (Status)
